$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cell values scraped for the cryptos list (coinranking.com) refresh.
$updates = @{
    'D2' = '29.253.20'
    'E2' = '  -0.39%  '
    'D3' = '1.863.82'
    'E3' = '  -1.10%  '
    'D4' = '1.003'
    'E4' = '  +0.28%  '
    'D5' = '242.59'
    'E5' = '  +0.10%  '
    'D6' = '0.7034'
    'E6' = '  -1.41%  '
    'D7' = '1.003'
    'E7' = '  +0.29%  '
    'D8' = '0.07813'
    'E8' = '  -3.52%  '
    'D9' = '0.3108'
    'E9' = '  -0.91%  '
    'D10' = '24.22'
    'E10' = '  -4.41%  '
    'B11' = 'TRON'
    'C11' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D11' = '0.08019'
    'E11' = '  -4.07%  '
    'B12' = 'WrappedEther'
    'C12' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
    'D12' = '1.913.08'
    'E12' = '  +1.77%  '
    'D13' = '5.169'
    'E13' = '  -1.58%  '
    'D14' = '92.96'
    'E14' = '  +0.87%  '
    'D15' = '0.6943'
    'E15' = '  -3.88%  '
    'D16' = '6.353'
    'E16' = '  +1.05%  '
    'D17' = '29.796.27'
    'E17' = '  +1.43%  '
    'D18' = '0.000008253'
    'E18' = '  -2.57%  '
    'B19' = 'WrappedliquidstakedEther2.0'
    'C19' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
    'D19' = '2.219.46'
    'E19' = '  +4.70%  '
    'B20' = 'BitcoinCash'
    'C20' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
    'D20' = '249.26'
    'E20' = '  +3.14%  '
    'D21' = '13.15'
    'E21' = '  -0.77%  '
    'E22' = '  +0.33%  '
    'D23' = '7.577'
    'E23' = '  -3.01%  '
    'D24' = '1.003'
    'E24' = '  +0.28%  '
    'D25' = '0.1547'
    'E25' = '  -2.95%  '
    'D26' = '8.967'
    'E26' = '  -1.28%  '
    'D27' = '159.86'
    'D28' = '18.63'
    'E28' = '  +0.19%  '
    'E29' = '  -0.56%  '
    'D30' = '4.268'
    'E30' = '  -2.07%  '
    'D31' = '4.272'
    'E31' = '  -3.57%  '
    'D32' = '1.223'
    'E32' = '  +0.28%  '
    'D33' = '0.05237'
    'E33' = '  -2.75%  '
    'D34' = '1.885'
    'E34' = '  -3.68%  '
    'D35' = '0.7440'
    'E35' = '  -1.18%  '
    'E36' = '  -2.04%  '
    'D37' = '2.715'
    'E37' = '  +0.72%  '
    'D38' = '0.01855'
    'E38' = '  -1.48%  '
    'D39' = '1.249.69'
    'E39' = '  -2.78%  '
    'D40' = '2.742'
    'E40' = '  -0.12%  '
    'D41' = '6.278'
    'E41' = '  -4.46%  '
    'D42' = '111.21'
    'E42' = '  +0.75%  '
    'D43' = '0.8963'
    'E43' = '  +0.33%  '
    'D44' = '71.58'
    'E44' = '  -2.78%  '
    'E45' = '  +0.30%  '
    'D46' = '0.00000000129'
    'E46' = '  +0.13%  '
    'D47' = '2.061.96'
    'E47' = '  +2.03%  '
    'B48' = 'RenderToken'
    'C48' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
    'D48' = '1.797'
    'E48' = '  -0.48%  '
    'B49' = 'Mantle'
    'C49' = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
    'D49' = '0.5199'
    'E49' = '  -0.30%  '
    'D51' = '1.014'
    'E51' = '  +1.50%  '
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so numeric-looking strings (e.g. "1.003") are not
    # re-interpreted by Excel as numbers, matching the original inline-string cells.
    $cell.NumberFormat = '@'
    $cell.Value = $updates[$cellRef]
}
